$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet 1: "USS Tester (2 Systems)"
# -------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Flip the two mismatched truth/sensed values on row 9 before trimming rows
$ws1.Cells.Item(9, 2).Value = 0   # B9: 1 -> 0
$ws1.Cells.Item(9, 3).Value = 0   # C9: 1 -> 0

# Drop rows 10-13 (rows shift up automatically)
$ws1.Range("A10:H13").EntireRow.Delete()

# Re-anchor the colour-scale conditional format to the new used range
$cf1 = $ws1.Range("H2:H9").FormatConditions.Item(1)
$cf1.ModifyAppliesToRange($ws1.Range("H2:H9"))

# Unfreeze panes / drop the split view
$ws1.Application.ActiveWindow.FreezePanes = $false

# -------------------------------------------------------------------------
# Sheet 2: "System 1 History"
# -------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(9, 2).Value = 0   # B9: 1 -> 0
$ws2.Cells.Item(9, 6).Value = 0   # F9: 1 -> 0

$ws2.Range("A10:L13").EntireRow.Delete()

$ws2.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($ws2.Range("D2:D10"))
$ws2.Range("E2:E10").FormatConditions.Item(1).ModifyAppliesToRange($ws2.Range("E2:E10"))
$ws2.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($ws2.Range("I2:I10"))
$ws2.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($ws2.Range("J2:J10"))
$ws2.Range("L2:L9").FormatConditions.Item(1).ModifyAppliesToRange($ws2.Range("L2:L9"))

# -------------------------------------------------------------------------
# Sheet 3: "System 2 History"
# -------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(7, 6).Value = 0    # F7: 1 -> 0
$ws3.Cells.Item(7, 13).Value = 0   # M7: 1 -> 0
$ws3.Cells.Item(8, 6).Value = 0    # F8: 1 -> 0
$ws3.Cells.Item(8, 13).Value = 0   # M8: 1 -> 0
$ws3.Cells.Item(9, 6).Value = 0    # F9: 1 -> 0
$ws3.Cells.Item(9, 13).Value = 0   # M9: 1 -> 0

$ws3.Range("A10:P13").EntireRow.Delete()

$ws3.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("D2:D10"))
$ws3.Range("E2:E10").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("E2:E10"))
$ws3.Range("F2:F10").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("F2:F10"))
$ws3.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("G2:G10"))
$ws3.Range("K2:K10").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("K2:K10"))
$ws3.Range("L2:L10").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("L2:L10"))
$ws3.Range("M2:M10").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("M2:M10"))
$ws3.Range("N2:N10").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("N2:N10"))
$ws3.Range("P2:P9").FormatConditions.Item(1).ModifyAppliesToRange($ws3.Range("P2:P9"))
